$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-range labels in column C (shared strings "8:20-8:25" / "8:25-8:30"
# become "8:35-8:40" / "8:40-8:45")
$ws.Range("C2").Value = "8:35-8:40"
$ws.Range("C3").Value = "8:40-8:45"

# Move the active selection from B12 to C14
$excel.Goto($ws.Range("C14"))
